# Update the Implementation Guide metadata (commit: "Atualizacao do Implementation Guide.")
#
# Changes applied:
#  1. Metadata!B2 (and Elements!R4, which shares the same text) - the
#     canonical URL host moves from "www.gabriellesantosleandro.com/molic-avc"
#     to "molic-avc.gabriellesantosleandro.com".
#  2. Metadata!B8 - the publication Date is bumped.
#  3. Elements!Z6 - the ValueSet URL host is rewritten the same way.
#  4. Elements column Z width shrinks to match the new (shorter) text's
#     best-fit width.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$newSdUrl = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/molicavc-genderidentity-extension"
$newVsUrl = "https://molic-avc.gabriellesantosleandro.com/ValueSet/molicavc-gender-valueset"

# 1. StructureDefinition canonical URL (Metadata sheet "URL" row).
$metadata.Range("B2").Value = $newSdUrl

# 2. Date row.
$metadata.Range("B8").Value = "2023-08-16T00:27:03-03:00"

# 3. The Elements grid repeats the StructureDefinition URL as the fixed
#    value of Extension.url (row 4, column R) - keep it in sync too.
$elements.Range("R4").Value = $newSdUrl

# 4. ValueSet URL referenced by the required binding.
$elements.Range("Z6").Value = $newVsUrl

# 5. Column Z (26) is best-fit sized to its content; after shortening the
#    text above, its stored width shrinks from 80.0234375 to 74.453125
#    characters. The COM layer only exposes ColumnWidth in coarse
#    (1/6-character) steps, so use the closest reachable value.
$elements.Columns.Item(26).ColumnWidth = 73.67
